# Updated symbol list (coin prices/volumes) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.252"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05701"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.298"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8099"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Value = "'0.1423"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07343"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03037"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03077"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09382"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.882"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001585"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04781"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005848"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006163"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005162"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.0009959"
$ws.Range("D21").Style = "Normal"
$ws.Range("D40").Value = "'0.03924"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006765"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1066"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Value = "'0.007493"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Value = "'0.1926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"
